# remove ExtendStems and roll back
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L/M column values (rows 7-11 and 17-21)
$ws.Range("L7").Value2 = 0.714
$ws.Range("M7").Value2 = 0.761

$ws.Range("L8").Value2 = 0.717
$ws.Range("M8").Value2 = 0.759

$ws.Range("L9").Value2 = 0.711
$ws.Range("M9").Value2 = 0.763

$ws.Range("L10").Value2 = 0.717
$ws.Range("M10").Value2 = 0.783

$ws.Range("L11").Value2 = 0.721
$ws.Range("M11").Value2 = 0.754

$ws.Range("L17").Value2 = 0.757
$ws.Range("M17").Value2 = 0.827

$ws.Range("L18").Value2 = 0.743
$ws.Range("M18").Value2 = 0.816

$ws.Range("L19").Value2 = 0.771
$ws.Range("M19").Value2 = 0.838

$ws.Range("L20").Value2 = 0.754
$ws.Range("M20").Value2 = 0.843

$ws.Range("L21").Value2 = 0.808
$ws.Range("M21").Value2 = 0.853

# Update the active cell selection
$ws.Range("I26").Select() | Out-Null
